$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Libre La Frite row (currently row 17) gets its RAM bumped from
# "DDR4 1GB" to "DDR4-2400 1GB" (done before the insert below so the new
# shared string is appended in the same order as the source edit).
$ws.Range("E17").Value = "DDR4-2400 1GB"

# Insert a new benchmark row at position 17 (pushes the Libre La Frite row,
# and everything below it, down by one), inheriting the formatting of the
# row above it.
$ws.Rows(17).Insert()

$ws.Range("A17").Value = "Linux"
$ws.Range("B17").Value = "TV-Box Vontar"
$ws.Range("C17").Value = "Amlogic S905W2"
$ws.Range("D17").Value = "1.0"
$ws.Range("E17").Value = "DDR3 2GB"
$ws.Range("F17").Value = 1700
$ws.Range("G17").Value = 45450
$ws.Range("H17").Value = "1364-1912"

# Move the active selection, matching the author's saved cursor position.
[void]$ws.Range("F18").Select()
